$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 content + alignment
$ws.Range("G5").Value = "100%, revisión en el sistema de prácticas"
$ws.Range("B5").Value = "Las observaciones obtenidas de la reunión del comité consultivo 2024"

# B2 merged cell: Career name
$ws.Range("B2").Value = "Ciencia de Datos e Inteligencia Artificial"

$ws.Range("C5").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B3").Select()
